# Apply cell updates for the cryptos price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.779.05"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.338.47"
$ws.Range("E3").Value = "  -0.85%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.669"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.43%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "237.65"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "72.74"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("E8").Value = "  -0.01%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.583"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.51%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0989"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.62%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "57.16"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "31.99"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.12%  "
$ws.Range("E13").Value = "  -0.10%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "2.686.94"
$ws.Range("E15").Value = "  -0.93%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "16.28"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -4.19%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.886"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("D18").Value = "2.327.40"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "43.609.60"
$ws.Range("E19").Value = "  -1.44%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000100"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.80"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.20%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "76.38"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "254.87"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +22.40%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -3.82%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.51%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.52"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E29").Value = "  -0.88%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "22.46"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "174.26"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.128"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("E33").Value = "  +1.19%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0746"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.53"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.01%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.12"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.92%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.84%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.33"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("E39").Value = "  -4.52%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0274"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +9.90%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.201"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.65%  "
$ws.Range("B43").Value = "BinanceUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.86"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.53"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.61%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "59.70"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +13.30%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.68"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.78%  "
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("E49").Value = "  -3.21%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "99.00"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  -2.33%  "
